$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New localization file being handed off: 721aa51a-e9be-4dc0-9833-32873f099577
# It is inserted alphabetically/chronologically between af6c4662... (row2) and
# 2aa438ea... (old row3) on every sheet, i.e. it becomes the new row 3 and all
# rows below shift down by one.
# ---------------------------------------------------------------------------

$newMdName   = "721aa51a-e9be-4dc0-9833-32873f099577.md"
$newZhXlf    = "721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.zh-cn.xlf"
$newDeXlf    = "721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.de-de.xlf"

$newCommit   = "c1a6c6fabf5a1f9c8b0a6a0a6b4a6f4a6b7c1a6c"

$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/$newCommit/e2e/$newMdName"
$zhXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZhXlf"
$deXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDeXlf"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(3).Insert()
$ws1.Range("A3").Value = $newMdName
$ws1.Range("B3").Value = "In Translation"
$ws1.Range("C3").Value = "In Translation"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = $newMdName
$ws2.Range("B3").Value = "In Translation"
$ws2.Range("C3").Value = $newZhXlf
$ws2.Range("D3").Value = "2016-03-03 06:46:52"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"
# The hand-back timestamp for the (now shifted down) 7cbb2cb4 row is new too.
$ws2.Range("D5").Value = "2016-03-03 06:49:09"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(3).Insert()
$ws3.Range("A3").Value = $newMdName
$ws3.Range("B3").Value = "In Translation"
$ws3.Range("C3").Value = $newDeXlf
$ws3.Range("D3").Value = "2016-03-03 06:47:11"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"
# The hand-back timestamp for the (now shifted down) 7cbb2cb4 row is new too.
$ws3.Range("D5").Value = "2016-03-03 06:49:20"

# ---------------------------------------------------------------------------
# Hyperlinks: row-insert does not shift hyperlink anchors, and Hyperlinks.Add
# does not replace an existing link on the same cell, so rebuild every link
# on every sheet from a clean slate once all the cell values are final.
# ---------------------------------------------------------------------------

$af6c4662Md  = "af6c4662-f8fd-4e34-957c-3654765d9d23.md"
$aa438eaMd   = "2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md"
$cbb2cb4Md   = "7cbb2cb4-f468-479a-965e-8eb53477f492.md"
$configName  = ".localization-config"

$af6c4662MdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/af6c4662-f8fd-4e34-957c-3654765d9d23.md"
$aa438eaMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/beb7b57e2f8af33a52a6a51a8e13cb8a2e7ad8f3/e2e/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.md"
$cbb2cb4MdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/beb7b57e2f8af33a52a6a51a8e13cb8a2e7ad8f3/e2e/7cbb2cb4-f468-479a-965e-8eb53477f492.md"
$configUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/beb7b57e2f8af33a52a6a51a8e13cb8a2e7ad8f3/.localization-config"

$af6c4662ZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e384b61201ce780c9dd60048116ca64bb0b41c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.zh-cn.xlf"
$aa438eaZhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2906c50a2f54f04353797ca69041d79cb2d0fe6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.3493345316d0d650da8b30231ef4f293442fe2f6.zh-cn.xlf"
$cbb2cb4ZhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2906c50a2f54f04353797ca69041d79cb2d0fe6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.zh-cn.xlf"

$af6c4662DeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32713c9abb62d7025c31384c79b02b15274b5191/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.de-de.xlf"
$aa438eaDeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11e35fea84d480bd46a97ac03cdaed6ca1fb2798/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2aa438ea-f3e3-428b-aff7-543eca8ba0f6.3493345316d0d650da8b30231ef4f293442fe2f6.de-de.xlf"
$cbb2cb4DeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11e35fea84d480bd46a97ac03cdaed6ca1fb2798/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.de-de.xlf"

# --- Sheet 1 hyperlinks (A column only) ---
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $af6c4662MdUrl, "", "", $af6c4662Md)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdUrl, "", "", $newMdName)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $aa438eaMdUrl, "", "", $aa438eaMd)
$ws1.Hyperlinks.Add($ws1.Range("A5"), $cbb2cb4MdUrl, "", "", $cbb2cb4Md)
$ws1.Hyperlinks.Add($ws1.Range("A6"), $configUrl, "", "", $configName)

# --- Sheet 2 (zh-cn) hyperlinks (A and C columns) ---
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $af6c4662MdUrl, "", "", $af6c4662Md)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $af6c4662ZhXlfUrl, "", "", "af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdUrl, "", "", $newMdName)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zhXlfUrl, "", "", $newZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $aa438eaMdUrl, "", "", $aa438eaMd)
$ws2.Hyperlinks.Add($ws2.Range("C4"), $aa438eaZhXlfUrl, "", "", "2aa438ea-f3e3-428b-aff7-543eca8ba0f6.3493345316d0d650da8b30231ef4f293442fe2f6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), $cbb2cb4MdUrl, "", "", $cbb2cb4Md)
$ws2.Hyperlinks.Add($ws2.Range("C5"), $cbb2cb4ZhXlfUrl, "", "", "7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), $configUrl, "", "", $configName)

# --- Sheet 3 (de-de) hyperlinks (A and C columns) ---
$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $af6c4662MdUrl, "", "", $af6c4662Md)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $af6c4662DeXlfUrl, "", "", "af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdUrl, "", "", $newMdName)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $deXlfUrl, "", "", $newDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $aa438eaMdUrl, "", "", $aa438eaMd)
$ws3.Hyperlinks.Add($ws3.Range("C4"), $aa438eaDeXlfUrl, "", "", "2aa438ea-f3e3-428b-aff7-543eca8ba0f6.3493345316d0d650da8b30231ef4f293442fe2f6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), $cbb2cb4MdUrl, "", "", $cbb2cb4Md)
$ws3.Hyperlinks.Add($ws3.Range("C5"), $cbb2cb4DeXlfUrl, "", "", "7cbb2cb4-f468-479a-965e-8eb53477f492.312cc3e5efb0e0c1fd7b54275018a861b3883926.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), $configUrl, "", "", $configName)

$wb.Save()
